# 27-A.xlsx — "Corrected many tasks in in2010402"
#
# Adds a second data series (F:G, mirroring the existing A:B pair) whose
# values are the existing A1:A40 numbers sorted in descending order, plus
# a MOD(.,3) helper column (G) built the same way column B was built, and
# three spot-check sum formulas in I1:I3 (mirroring D1:D3). Finally the F
# column is sorted (descending) via the UI Sort feature, which is what
# stamps the new <sortState> or keeps it consistent with the data already
# being in descending order, and the selection is left on I2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F1:F40 = A1:A40 reversed (i.e. the same 40 numbers, descending).
$descending = @(951,919,898,894,880,874,850,838,832,775,724,694,691,667,617,566,523,516,513,503,487,432,406,376,375,364,323,308,286,273,268,246,217,215,206,203,193,136,56,25)

for ($i = 0; $i -lt $descending.Length; $i++) {
    $ws.Cells.Item($i + 1, 6).Value = $descending[$i]
}

# G1:G40 = MOD(F,3), entered the same way B1:B40 was (single formula in
# row 1, then a fill-down that Excel stores as one shared formula).
$ws.Range("G1").Formula = "=MOD(F1,3)"
$ws.Range("G2:G40").Formula = "=MOD(F2,3)"

# I1:I3 = a few spot-check sums over the new F column (mirrors D1:D3).
$ws.Range("I1").Formula = "=F1+F2+F15"
$ws.Range("I2").Formula = "=F2+F3+F5"
$ws.Range("I3").Formula = "=F1+F4+F18"

# Sort F1:F40 descending (matches the data, and updates <sortState>).
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("F1"), 0, 2)
$ws.Sort.SetRange($ws.Range("F1:F40"))
$ws.Sort.Header = 0
$ws.Sort.Apply()

# Leave the selection on I2, as in the final workbook.
$null = $ws.Range("I2").Select()
